$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("marker_info")

# Insert a new column at H (shifts old H/I/J -> I/J/K)
$ws.Columns.Item(8).Insert()

# Set header for the new column
$ws.Range("H1").Value = "alternate ID3"

# Populate the two new values that were filled in for rows 13 and 23
$ws.Range("H13").Value = "Scaffold79929e:640165"
$ws.Range("H23").Value = "Scaffold79929e:670329"

# Match the column width used for the new column (ColumnWidth 40.67 -> stored width 41.5)
$ws.Columns.Item(8).ColumnWidth = 40.67

# Re-point the conditional formatting that used to cover H2:I36 to the shifted I2:J36 range,
# including the formula's relative cell reference
$fc = $ws.Range("H2:I36").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("I2:J36"))
$fc.Formula1 = '=NOT(ISERROR(SEARCH("FALSE",I2)))'

# Refresh the sort range so it covers the new column too
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("D2:D36"))
$ws.Sort.SetRange($ws.Range("A2:K36"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Update the selection to match the recorded view state
$ws.Range("H25").Select()
